# HERVAL.xlsx automatic update
#
# 1. Rename "Paineis DARQ" -> "PAINEIS DARQ"
# 2. Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# 3. Delete sheet "Desarquivamentos Pendentes" (no longer needed)

$wb = $excel.ActiveWorkbook

$wb.Worksheets("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

$wb.Worksheets("Desarquivamentos Pendentes").Delete()
